$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "69.261.03"
$ws.Range("E2").Value = "  -0.89%  "
$ws.Range("D3").Value = "3.514.25"
$ws.Range("E3").Value = "  -1.96%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").Value = "'574.05"
$ws.Range("E5").Value = "  -1.00%  "
$ws.Range("D6").Value = "'185.29"
$ws.Range("E6").Value = "  -2.89%  "
$ws.Range("D7").Value = "3.501.73"
$ws.Range("E7").Value = "  -2.15%  "
$ws.Range("D8").Value = "'0.614"
$ws.Range("E8").Value = "  -2.80%  "
$ws.Range("E9").Value = "  +0.08%  "
$ws.Range("D10").Value = "'0.186"
$ws.Range("E10").Value = "  +3.26%  "
$ws.Range("E11").Value = "  -1.94%  "
$ws.Range("D12").Value = "'54.35"
$ws.Range("E12").Value = "  -2.51%  "
$ws.Range("E13").Value = "  -1.21%  "
$ws.Range("D14").Value = "'9.46"
$ws.Range("E14").Value = "  -1.84%  "
$ws.Range("D15").Value = "4.069.85"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "'19.37"
$ws.Range("E16").Value = "  -2.44%  "
$ws.Range("B17").Value = "WrappedBTC"
$ws.Range("C17").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D17").Value = "69.102.16"
$ws.Range("E17").Value = "  -1.01%  "
$ws.Range("B18").Value = "WrappedEther"
$ws.Range("C18").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D18").Value = "3.503.34"
$ws.Range("E18").Value = "  -2.16%  "
$ws.Range("D19").Value = "'12.35"
$ws.Range("E19").Value = "  -2.44%  "
$ws.Range("E20").Value = "  -1.05%  "
$ws.Range("D21").Value = "'546.07"
$ws.Range("E21").Value = "  +15.10%  "
$ws.Range("E22").Value = "  -2.77%  "
$ws.Range("D23").Value = "'18.45"
$ws.Range("E23").Value = "  -5.32%  "
$ws.Range("D24").Value = "'4.94"
$ws.Range("E24").Value = "  -1.73%  "
$ws.Range("D25").Value = "'4.44"
$ws.Range("E25").Value = "  +1.53%  "
$ws.Range("D26").Value = "'94.10"
$ws.Range("E26").Value = "  -1.79%  "
$ws.Range("D27").Value = "'11.19"
$ws.Range("E27").Value = "  +1.45%  "
$ws.Range("D28").Value = "'2.94"
$ws.Range("E28").Value = "  -2.07%  "
$ws.Range("D29").Value = "'9.14"
$ws.Range("E29").Value = "  -1.54%  "
$ws.Range("D30").Value = "'31.75"
$ws.Range("E30").Value = "  -1.44%  "
$ws.Range("D31").Value = "'7.31"
$ws.Range("E31").Value = "  -4.22%  "
$ws.Range("D32").Value = "'12.72"
$ws.Range("E32").Value = "  +4.30%  "
$ws.Range("D33").Value = "'64.72"
$ws.Range("E33").Value = "  -2.46%  "
$ws.Range("E34").Value = "  -3.93%  "
$ws.Range("D35").Value = "'548.56"
$ws.Range("E35").Value = "  -7.09%  "
$ws.Range("B36").Value = "TheGraph"
$ws.Range("C36").Value = "https://coinranking.com/coin/qhd1biQ7M+thegraph-grt"
$ws.Range("D36").Value = "'0.406"
$ws.Range("E36").Value = "  +2.92%  "
$ws.Range("D37").Value = "'38.16"
$ws.Range("E37").Value = "  -2.13%  "
$ws.Range("B38").Value = "Fetch.AI"
$ws.Range("C38").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D38").Value = "'3.08"
$ws.Range("E38").Value = "  +7.86%  "
$ws.Range("B39").Value = "Dai"
$ws.Range("C39").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D39").Value = "'1.00"
$ws.Range("E39").Value = "  -0.06%  "
$ws.Range("D40").Value = "0.0₃0768"
$ws.Range("E40").Value = "  -3.90%  "
$ws.Range("B41").Value = "dogwifhat"
$ws.Range("C41").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D41").Value = "'3.11"
$ws.Range("E41").Value = "  -2.35%  "
$ws.Range("B42").Value = "Stacks"
$ws.Range("C42").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D42").Value = "'3.38"
$ws.Range("E42").Value = "  -2.39%  "
$ws.Range("E43").Value = "  -2.88%  "
$ws.Range("D44").Value = "3.310.71"
$ws.Range("E44").Value = "  +2.74%  "
$ws.Range("E45").Value = "  -2.64%  "
$ws.Range("D46").Value = "'0.0445"
$ws.Range("E46").Value = "  +0.71%  "
$ws.Range("E47").Value = "  +2.86%  "
$ws.Range("D48").Value = "'0.135"
$ws.Range("E48").Value = "  -2.22%  "
$ws.Range("D49").Value = "'8.93"
$ws.Range("E49").Value = "  -5.22%  "
$ws.Range("D50").Value = "'0.997"
$ws.Range("E50").Value = "  -0.32%  "
$ws.Range("D51").Value = "'137.43"
$ws.Range("E51").Value = "  +2.93%  "
